# "added data for choropleth"
#
# Semantic changes applied:
#  1. Rename the first worksheet from
#       "2016 State Manufacturing vs Che" -> "2016 State Checks vs Manufactur"
#  2. Make the first worksheet ("2016 State Checks vs Manufactur") the active /
#     selected tab (it previously was the second sheet, "Estimated Sales vs
#     Manufacturin", that was active).
#  3. Update the remembered selection on each sheet:
#       - Sheet 1 ("2016 State Checks vs Manufactur"): selected cell stays A40
#         (selection unchanged, but the sheet's scroll position resets to the
#         top of the frozen pane, A2, once it becomes the active sheet).
#       - Sheet 2 ("Estimated Sales vs Manufacturin"): selected cell moves
#         from H8 to J12.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename the first sheet.
$ws1.Name = "2016 State Checks vs Manufactur"

# 2/3. Update the selection remembered on sheet 2, then sheet 1, finally
# leaving sheet 1 as the active (selected) tab.
$ws2.Activate()
$null = $ws2.Range("J12").Select()

$ws1.Activate()
$null = $ws1.Range("A40").Select()
